$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" '26.382.23'
Set-TextValue "E2" '  -3.34%  '
Set-TextValue "D3" '1.800.06'
Set-TextValue "E3" '  -2.94%  '
Set-TextValue "D4" '1.007'
Set-TextValue "E4" '  +0.49%  '
Set-TextValue "D5" '1.008'
Set-TextValue "E5" '  +0.65%  '
Set-TextValue "D6" '307.96'
Set-TextValue "E6" '  -1.98%  '
Set-TextValue "D7" '0.4530'
Set-TextValue "E7" '  -1.60%  '
Set-TextValue "D8" '0.3638'
Set-TextValue "E8" '  -1.92%  '
Set-TextValue "D9" '0.07078'
Set-TextValue "E9" '  -3.02%  '
Set-TextValue "D10" '0.8687'
Set-TextValue "E10" '  -2.16%  '
Set-TextValue "D11" '0.07759'
Set-TextValue "E11" '  -0.79%  '
Set-TextValue "D12" '19.22'
Set-TextValue "E12" '  -4.58%  '
Set-TextValue "D13" '1.806.55'
Set-TextValue "E13" '  -1.94%  '
Set-TextValue "D14" '5.239'
Set-TextValue "E14" '  -2.75%  '
Set-TextValue "D15" '6.311'
Set-TextValue "E15" '  -3.37%  '
Set-TextValue "D16" '85.67'
Set-TextValue "E16" '  -6.30%  '
Set-TextValue "E17" '  +0.66%  '
Set-TextValue "D18" '0.000008526'
Set-TextValue "E18" '  -4.55%  '
Set-TextValue "D19" '1.007'
Set-TextValue "E19" '  +0.55%  '
Set-TextValue "D20" '26.443.72'
Set-TextValue "E20" '  -3.14%  '
Set-TextValue "E21" '  -4.02%  '
Set-TextValue "D22" '4.952'
Set-TextValue "E22" '  -3.17%  '
Set-TextValue "D23" '10.36'
Set-TextValue "E23" '  -1.83%  '
Set-TextValue "D24" '1.967'
Set-TextValue "E24" '  +2.26%  '
Set-TextValue "D25" '150.61'
Set-TextValue "E25" '  -0.97%  '
Set-TextValue "D26" '17.84'
Set-TextValue "E26" '  -3.47%  '
Set-TextValue "D27" '1.975'
Set-TextValue "E27" '  -4.19%  '
Set-TextValue "D28" '112.54'
Set-TextValue "E28" '  -3.01%  '
Set-TextValue "D29" '4.849'
Set-TextValue "E29" '  -4.38%  '
Set-TextValue "D30" '0.08640'
Set-TextValue "E30" '  -2.11%  '
Set-TextValue "D31" '3.035'
Set-TextValue "E31" '  -1.38%  '
Set-TextValue "D32" '0.7248'
Set-TextValue "E32" '  -6.29%  '
Set-TextValue "D33" '4.423'
Set-TextValue "E33" '  -1.86%  '
Set-TextValue "D34" '1.104'
Set-TextValue "E34" '  -5.70%  '
Set-TextValue "D35" '1.005'
Set-TextValue "E35" '  +0.44%  '
Set-TextValue "D36" '2.508'
Set-TextValue "E36" '  -9.10%  '
Set-TextValue "D37" '1.074'
Set-TextValue "E37" '  -0.65%  '
Set-TextValue "D39" '2.871'
Set-TextValue "E39" '  -2.78%  '
Set-TextValue "D40" '0.05048'
Set-TextValue "E40" '  -3.89%  '
Set-TextValue "D41" '6.930'
Set-TextValue "E41" '  -1.86%  '
Set-TextValue "D42" '0.4898'
Set-TextValue "E42" '  -4.42%  '
Set-TextValue "D43" '0.1563'
Set-TextValue "D44" '8.078'
Set-TextValue "E44" '  -3.93%  '
Set-TextValue "D45" '1.009'
Set-TextValue "E45" '  +0.80%  '
Set-TextValue "D46" '0.4579'
Set-TextValue "E46" '  -4.57%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D47" '101.13'
Set-TextValue "E47" '  -1.18%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D48" '9.868'
Set-TextValue "E48" '  -4.54%  '
Set-TextValue "D49" '1.574'
Set-TextValue "E49" '  -4.24%  '
Set-TextValue "E50" '  -3.80%  '
Set-TextValue "D51" '63.31'
Set-TextValue "E51" '  -3.63%  '
